$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuilt "Estado de Cuenta" detail table (rows 16-48): previous statements
# removed and new ones added, data grouped per worker across all periods.
$rows = @(
    @{ Row=16; B='CC'; C='1047487411'; D='LUNA LUZ TRESPALACIO RODRIGUEZ'; E='2205'; F=27861; G=877802 },
    @{ Row=17; B='CC'; C='1047487411'; D='LUNA LUZ TRESPALACIO RODRIGUEZ'; E='2204'; F=36341; G=877802 },
    @{ Row=18; B='CC'; C='1047487411'; D='LUNA LUZ TRESPALACIO RODRIGUEZ'; E='2203'; F=36341; G=877802 },
    @{ Row=19; B='CC'; C='1047487411'; D='LUNA LUZ TRESPALACIO RODRIGUEZ'; E='2202'; F=36341; G=877802 },
    @{ Row=20; B='CC'; C='1047487411'; D='LUNA LUZ TRESPALACIO RODRIGUEZ'; E='2201'; F=36341; G=877802 },
    @{ Row=21; B='CC'; C='1047487411'; D='LUNA LUZ TRESPALACIO RODRIGUEZ'; E='2112'; F=36341; G=877802 },
    @{ Row=22; B='CC'; C='1047487411'; D='LUNA LUZ TRESPALACIO RODRIGUEZ'; E='2111'; F=36341; G=877802 },
    @{ Row=23; B='CC'; C='1047487411'; D='LUNA LUZ TRESPALACIO RODRIGUEZ'; E='2110'; F=36341; G=877802 },
    @{ Row=24; B='CC'; C='1044920194'; D='VIKI PAOLA PAJARO ZAMBRANO'; E='2205'; F=27861; G=908526 },
    @{ Row=25; B='CC'; C='1044920194'; D='VIKI PAOLA PAJARO ZAMBRANO'; E='2204'; F=36341; G=908526 },
    @{ Row=26; B='CC'; C='1044920194'; D='VIKI PAOLA PAJARO ZAMBRANO'; E='2203'; F=36341; G=908526 },
    @{ Row=27; B='CC'; C='1044920194'; D='VIKI PAOLA PAJARO ZAMBRANO'; E='2202'; F=36341; G=908526 },
    @{ Row=28; B='CC'; C='1044920194'; D='VIKI PAOLA PAJARO ZAMBRANO'; E='2201'; F=36341; G=908526 },
    @{ Row=29; B='CC'; C='1044920194'; D='VIKI PAOLA PAJARO ZAMBRANO'; E='2112'; F=36341; G=908526 },
    @{ Row=30; B='CC'; C='1044920194'; D='VIKI PAOLA PAJARO ZAMBRANO'; E='2111'; F=36341; G=908526 },
    @{ Row=31; B='CC'; C='1044920194'; D='VIKI PAOLA PAJARO ZAMBRANO'; E='2110'; F=36341; G=908526 },
    @{ Row=32; B='CC'; C='1143351410'; D='YOHANA DE JESUS QUINTANA OSORIO'; E='2110'; F=36341; G=908526 },
    @{ Row=33; B='CC'; C='1047504169'; D='MARYSOL TRESPALACIOS RODRIGUEZ'; E='2205'; F=27861; G=877803 },
    @{ Row=34; B='CC'; C='1047504169'; D='MARYSOL TRESPALACIOS RODRIGUEZ'; E='2204'; F=36341; G=877803 },
    @{ Row=35; B='CC'; C='1047504169'; D='MARYSOL TRESPALACIOS RODRIGUEZ'; E='2203'; F=36341; G=877803 },
    @{ Row=36; B='CC'; C='1047504169'; D='MARYSOL TRESPALACIOS RODRIGUEZ'; E='2202'; F=36341; G=877803 },
    @{ Row=37; B='CC'; C='1047504169'; D='MARYSOL TRESPALACIOS RODRIGUEZ'; E='2201'; F=36341; G=877803 },
    @{ Row=38; B='CC'; C='1047504169'; D='MARYSOL TRESPALACIOS RODRIGUEZ'; E='2112'; F=36341; G=877803 },
    @{ Row=39; B='CC'; C='1047504169'; D='MARYSOL TRESPALACIOS RODRIGUEZ'; E='2111'; F=36341; G=877803 },
    @{ Row=40; B='CC'; C='1047504169'; D='MARYSOL TRESPALACIOS RODRIGUEZ'; E='2110'; F=36341; G=877803 },
    @{ Row=41; B='CC'; C='1047498678'; D='LUIS EMIRO MONTALVO HERRERA'; E='2205'; F=27861; G=908526 },
    @{ Row=42; B='CC'; C='1047498678'; D='LUIS EMIRO MONTALVO HERRERA'; E='2204'; F=36341; G=908526 },
    @{ Row=43; B='CC'; C='1047498678'; D='LUIS EMIRO MONTALVO HERRERA'; E='2203'; F=36341; G=908526 },
    @{ Row=44; B='CC'; C='1047498678'; D='LUIS EMIRO MONTALVO HERRERA'; E='2202'; F=36341; G=908526 },
    @{ Row=45; B='CC'; C='1047498678'; D='LUIS EMIRO MONTALVO HERRERA'; E='2201'; F=36341; G=908526 },
    @{ Row=46; B='CC'; C='1047498678'; D='LUIS EMIRO MONTALVO HERRERA'; E='2112'; F=36341; G=908526 },
    @{ Row=47; B='CC'; C='1047498678'; D='LUIS EMIRO MONTALVO HERRERA'; E='2111'; F=36341; G=908526 },
    @{ Row=48; B='CC'; C='1047498678'; D='LUIS EMIRO MONTALVO HERRERA'; E='2110'; F=36341; G=908526 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
}